$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows at the bottom (351:355), inheriting formatting from the
# row immediately above (350), which still carries the legacy highlight style.
$ws.Rows("351:355").Insert(-4121, 0)

# New Lotomania draw results
$newData = @(
  @(2857,4,12,13,20,22,26,30,34,36,52,69,74,77,80,82,87,92,95,96,99),
  @(2858,1,3,6,8,9,20,24,25,43,44,52,58,61,66,67,77,88,89,90,93),
  @(2859,2,11,12,14,17,19,20,28,31,32,41,51,55,68,73,75,77,90,98,99),
  @(2860,0,1,5,6,9,12,14,22,25,29,31,34,39,47,66,67,72,79,80,98),
  @(2861,4,13,17,22,29,32,36,46,50,57,66,71,75,77,81,86,89,90,95,99)
)

$r = 351
foreach ($row in $newData) {
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r++
}

# The previously highlighted block (the "last 6 draws") loses its highlight
# now that the new draws have taken its place.
$ws.Range("A345:U350").Style = "Normal"

# Reflect the new selection/active cell on the freshly added block.
$null = $ws.Range("B351:U355").Select()
